# Generate Report for Handoff
# Updates the "localization-status" workbook: the two files that were
# previously blocked on handoff (fc9bd0dc-...md and .localization-config)
# now show a successful handoff, and a brand-new file
# (ffff01b68f77-...md) shows up as also ready for handoff.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$oldMdFile  = "fc9bd0dc-745f-45f5-8096-bf9b461f5a24.md"
$newMdFile  = "366d26ba-c56f-42f6-8320-c2b4558e46c9.md"
$newMdFile2 = "ffff01b68f77-91a0-458a-b92c-05df0b7578d0.md"
$cfgFile    = ".localization-config"

$zhXlf = "366d26ba-c56f-42f6-8320-c2b4558e46c9.f71aeb345eab118bfcd9e86a546e64b28d90b546.zh-cn.xlf"
$deXlf = "366d26ba-c56f-42f6-8320-c2b4558e46c9.f71aeb345eab118bfcd9e86a546e64b28d90b546.de-de.xlf"

$zhHandoffDt = "2016-01-27 08:22:47"
$deHandoffDt = "2016-01-27 08:23:00"
$epoch       = "0001-01-01 00:00:00"

$repoBase = "https://github.com/OpenLocalizationTest/oltest/blob"
$mdCommit = "6bb0f631995bc33ad54c39ced02eac6a756347b7"
$cfgCommit = "00b018a60bf10950d7b78faa5eeb9a3863907515"
$xlfCommit = "f71aeb345eab118bfcd9e86a546e64b28d90b546"

function Set-ReportRows($ws, [bool]$isDetail, [string]$xlfName, [string]$handoffDt) {
    # Insert a new row 3 (pushes the old row 3 -> row 4), then fill in the
    # three data rows with their final contents.
    $ws.Rows.Item(3).Insert()

    # Clear out any stale hyperlinks (their refs/targets no longer line up
    # after the insert); we'll recreate every hyperlink below.
    $ws.Hyperlinks.Delete()

    if (-not $isDetail) {
        # ---- Overview sheet: columns A (File Name), B (zh-cn), C (de-de) ----
        $ws.Range("A2").Value = $newMdFile
        $ws.Range("B2").Value = "Ready for handoff"
        $ws.Range("C2").Value = "Ready for handoff"

        $ws.Range("A3").Value = $newMdFile2
        $ws.Range("B3").Value = "Ready for handoff"
        $ws.Range("C3").Value = "Ready for handoff"

        $ws.Range("A4").Value = $cfgFile
        $ws.Range("B4").Value = "Not to be localized"
        $ws.Range("C4").Value = "Not to be localized"

        $ws.Hyperlinks.Add($ws.Range("A2"), "$repoBase/$mdCommit/e2e/$newMdFile", $null, $null, $newMdFile)
        $ws.Hyperlinks.Add($ws.Range("A3"), "$repoBase/$mdCommit/e2e/$newMdFile2", $null, $null, $newMdFile2)
        $ws.Hyperlinks.Add($ws.Range("A4"), "$repoBase/$cfgCommit/$cfgFile", $null, $null, $cfgFile)

        # Match the hyperlink look already used elsewhere in this sheet
        # (underlined, cornflower-blue text) instead of Excel's default
        # theme hyperlink color.
        $ws.Range("A2:A4").Font.Underline = 2
        $ws.Range("A2:A4").Font.Color = 15570276
    } else {
        # ---- Detail sheets (zh-cn / de-de): A Source File Name, B Status,
        #      C Latest Handoff File, D Latest Handoff Datetime,
        #      G Latest Handback DateTime, H Handoff Reason ----
        $ws.Range("A2").Value = $newMdFile
        $ws.Range("B2").Value = "Ready for handoff"
        $ws.Range("C2").Value = $xlfName
        $ws.Range("D2").Value = $handoffDt
        $ws.Range("G2").Value = $epoch
        $ws.Range("H2").Value = "Include"

        $ws.Range("A3").Value = $newMdFile2
        $ws.Range("B3").Value = "Ready for handoff"
        $ws.Range("C3").Value = $xlfName
        $ws.Range("D3").Value = $handoffDt
        $ws.Range("G3").Value = $epoch
        $ws.Range("H3").Value = "Include"

        $ws.Range("A4").Value = $cfgFile
        $ws.Range("B4").Value = "Not to be localized"
        $ws.Range("D4").Value = $epoch
        $ws.Range("G4").Value = $epoch
        $ws.Range("H4").Value = "Ignored"

        $ws.Range("D2:D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

        $ws.Hyperlinks.Add($ws.Range("A2"), "$repoBase/$mdCommit/e2e/$newMdFile", $null, $null, $newMdFile)
        $ws.Hyperlinks.Add($ws.Range("C2"), "$repoBase/$xlfCommit/handoff/$xlfName", $null, $null, $xlfName)
        $ws.Hyperlinks.Add($ws.Range("A3"), "$repoBase/$mdCommit/e2e/$newMdFile2", $null, $null, $newMdFile2)
        $ws.Hyperlinks.Add($ws.Range("C3"), "$repoBase/$xlfCommit/handoff/$xlfName", $null, $null, $xlfName)
        $ws.Hyperlinks.Add($ws.Range("A4"), "$repoBase/$cfgCommit/$cfgFile", $null, $null, $cfgFile)

        # Match the hyperlink look already used elsewhere in this sheet
        # (underlined, cornflower-blue text) instead of Excel's default
        # theme hyperlink color.
        $ws.Range("A2:A4").Font.Underline = 2
        $ws.Range("A2:A4").Font.Color = 15570276
        $ws.Range("C2:C3").Font.Underline = 2
        $ws.Range("C2:C3").Font.Color = 15570276
    }
}

Set-ReportRows $ws1 $false "" ""
Set-ReportRows $ws2 $true $zhXlf $zhHandoffDt
Set-ReportRows $ws3 $true $deXlf $deHandoffDt
